# SOPORTE U LIDIMA.xlsx - refactor code and structure files
#
# Update the price-list date (A1) and the unit price (D27:D29) on the
# "Hoja1" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date stamp at the top of the sheet: 24/04/2024 -> 24/05/2024
$ws.Range("A1").Value = 45436

# Unit price for the three "U" support items went from 94.30 to 203.50
$ws.Range("D27").Value = 203.5
$ws.Range("D28").Value = 203.5
$ws.Range("D29").Value = 203.5
